# vm_pu.xlsx ("case with 380 kV done")
# Re-run of the power-flow case with the slack-bus voltage setpoint lowered
# from 1.05 pu to 1.02 pu (column B). All other per-bus voltage-magnitude
# results (columns C:F and I:N, data rows 2-25) are refreshed with the
# recomputed values from that new solve. Column G (fixed at 1) and the blank
# column H are left untouched, as is the bus-index column A and the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033366029227542
$ws.Cells.Item(2, 4).Value = 1.025507051647275
$ws.Cells.Item(2, 5).Value = 1.041296853725589
$ws.Cells.Item(2, 6).Value = 1.049333393264743
$ws.Cells.Item(2, 9).Value = 1.02804731065926
$ws.Cells.Item(2, 10).Value = 1.038491151223592
$ws.Cells.Item(2, 11).Value = 1.028332420894444
$ws.Cells.Item(2, 12).Value = 1.044076768678686
$ws.Cells.Item(2, 13).Value = 1.052090743175996
$ws.Cells.Item(2, 14).Value = 1.016359438818662

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.035456500123056
$ws.Cells.Item(3, 4).Value = 1.026019325364455
$ws.Cells.Item(3, 5).Value = 1.043210675826196
$ws.Cells.Item(3, 6).Value = 1.051452051245787
$ws.Cells.Item(3, 9).Value = 1.028114196249099
$ws.Cells.Item(3, 10).Value = 1.040218527726439
$ws.Cells.Item(3, 11).Value = 1.028652813399939
$ws.Cells.Item(3, 12).Value = 1.045798295184488
$ws.Cells.Item(3, 13).Value = 1.054018255724803
$ws.Cells.Item(3, 14).Value = 1.016972109842518

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.036804338530668
$ws.Cells.Item(4, 4).Value = 1.026344013162027
$ws.Cells.Item(4, 5).Value = 1.044444645201334
$ws.Cells.Item(4, 6).Value = 1.052818661432289
$ws.Cells.Item(4, 9).Value = 1.028152502439901
$ws.Cells.Item(4, 10).Value = 1.041331346126495
$ws.Cells.Item(4, 11).Value = 1.028852644648525
$ws.Cells.Item(4, 12).Value = 1.046907440178718
$ws.Cells.Item(4, 13).Value = 1.055260824523872
$ws.Cells.Item(4, 14).Value = 1.017365833892181

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.037369841693726
$ws.Cells.Item(5, 4).Value = 1.026478881439717
$ws.Cells.Item(5, 5).Value = 1.044962378131683
$ws.Cells.Item(5, 6).Value = 1.053392182731227
$ws.Cells.Item(5, 9).Value = 1.028167412619662
$ws.Cells.Item(5, 10).Value = 1.041798024285772
$ws.Cells.Item(5, 11).Value = 1.028934856412294
$ws.Cells.Item(5, 12).Value = 1.047372600427122
$ws.Cells.Item(5, 13).Value = 1.055782111355798
$ws.Cells.Item(5, 14).Value = 1.017530713907021

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03746472686782
$ws.Cells.Item(6, 4).Value = 1.026501430639359
$ws.Cells.Item(6, 5).Value = 1.045049248248012
$ws.Cells.Item(6, 6).Value = 1.053488421571141
$ws.Cells.Item(6, 9).Value = 1.028169846008701
$ws.Cells.Item(6, 10).Value = 1.04187631483863
$ws.Cells.Item(6, 11).Value = 1.028948554567338
$ws.Cells.Item(6, 12).Value = 1.047450637651736
$ws.Cells.Item(6, 13).Value = 1.055869574615349
$ws.Cells.Item(6, 14).Value = 1.01755836060787

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.036811899204402
$ws.Cells.Item(7, 4).Value = 1.026345821692126
$ws.Cells.Item(7, 5).Value = 1.0444515671777
$ws.Cells.Item(7, 6).Value = 1.052826328749198
$ws.Cells.Item(7, 9).Value = 1.028152706365016
$ws.Cells.Item(7, 10).Value = 1.0413375863911
$ws.Cells.Item(7, 11).Value = 1.028853750235325
$ws.Cells.Item(7, 12).Value = 1.046913660057774
$ws.Cells.Item(7, 13).Value = 1.055267794227486
$ws.Cells.Item(7, 14).Value = 1.017368039536576

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.034073532951003
$ws.Cells.Item(8, 4).Value = 1.02568157912909
$ws.Cells.Item(8, 5).Value = 1.041944565930772
$ws.Cells.Item(8, 6).Value = 1.050050311958026
$ws.Cells.Item(8, 9).Value = 1.02807094286734
$ws.Cells.Item(8, 10).Value = 1.039075957746021
$ws.Cells.Item(8, 11).Value = 1.028442245416833
$ws.Cells.Item(8, 12).Value = 1.04465957474345
$ws.Cells.Item(8, 13).Value = 1.052743136263629
$ws.Cells.Item(8, 14).Value = 1.016567061135038

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029209801128964
$ws.Cells.Item(9, 4).Value = 1.024459342914767
$ws.Cells.Item(9, 5).Value = 1.037492025308604
$ws.Cells.Item(9, 6).Value = 1.045124365226372
$ws.Cells.Item(9, 9).Value = 1.027888906352644
$ws.Cells.Item(9, 10).Value = 1.035051958522594
$ws.Cells.Item(9, 11).Value = 1.027660033183786
$ws.Cells.Item(9, 12).Value = 1.040649743577527
$ws.Cells.Item(9, 13).Value = 1.048257462646569
$ws.Cells.Item(9, 14).Value = 1.015134458890293

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025939617291363
$ws.Cells.Item(10, 4).Value = 1.023610008354521
$ws.Cells.Item(10, 5).Value = 1.034498556500015
$ws.Cells.Item(10, 6).Value = 1.041815550827719
$ws.Cells.Item(10, 9).Value = 1.027742192411911
$ws.Cells.Item(10, 10).Value = 1.032341682945349
$ws.Cells.Item(10, 11).Value = 1.027100466945094
$ws.Cells.Item(10, 12).Value = 1.037949528146233
$ws.Cells.Item(10, 13).Value = 1.045240498997322
$ws.Cells.Item(10, 14).Value = 1.014164613834184

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024516597594767
$ws.Cells.Item(11, 4).Value = 1.02323410637176
$ws.Cells.Item(11, 5).Value = 1.033196024245224
$ws.Cells.Item(11, 6).Value = 1.040376494231811
$ws.Cells.Item(11, 9).Value = 1.027672681232815
$ws.Cells.Item(11, 10).Value = 1.031161200832587
$ws.Cells.Item(11, 11).Value = 1.026849189472056
$ws.Cells.Item(11, 12).Value = 1.036773553826061
$ws.Cells.Item(11, 13).Value = 1.043927448861385
$ws.Cells.Item(11, 14).Value = 1.013741033680938

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.023986934918851
$ws.Cells.Item(12, 4).Value = 1.023093261465912
$ws.Cells.Item(12, 5).Value = 1.032711220915079
$ws.Cells.Item(12, 6).Value = 1.039840979805725
$ws.Cells.Item(12, 9).Value = 1.027645965106361
$ws.Cells.Item(12, 10).Value = 1.030721647574989
$ws.Cells.Item(12, 11).Value = 1.026754508480939
$ws.Cells.Item(12, 12).Value = 1.036335698255753
$ws.Cells.Item(12, 13).Value = 1.043438686032643
$ws.Cells.Item(12, 14).Value = 1.013583141053325

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024100599274606
$ws.Cells.Item(13, 4).Value = 1.02312352824074
$ws.Cells.Item(13, 5).Value = 1.032815257965506
$ws.Cells.Item(13, 6).Value = 1.039955894605354
$ws.Cells.Item(13, 9).Value = 1.027651736342748
$ws.Cells.Item(13, 10).Value = 1.030815982143403
$ws.Cells.Item(13, 11).Value = 1.02677487869296
$ws.Cells.Item(13, 12).Value = 1.036429667602273
$ws.Cells.Item(13, 13).Value = 1.043543574808476
$ws.Cells.Item(13, 14).Value = 1.013617034894785

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024472837959735
$ws.Cells.Item(14, 4).Value = 1.023222488916466
$ws.Cells.Item(14, 5).Value = 1.033155970530062
$ws.Cells.Item(14, 6).Value = 1.040332248738053
$ws.Cells.Item(14, 9).Value = 1.027670491146611
$ws.Cells.Item(14, 10).Value = 1.031124889223288
$ws.Cells.Item(14, 11).Value = 1.026841390534216
$ws.Cells.Item(14, 12).Value = 1.036737382070019
$ws.Cells.Item(14, 13).Value = 1.043887068933161
$ws.Cells.Item(14, 14).Value = 1.013727993633358

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024702040976964
$ws.Cells.Item(15, 4).Value = 1.023283300548361
$ws.Cells.Item(15, 5).Value = 1.033365763308621
$ws.Cells.Item(15, 6).Value = 1.040564001381583
$ws.Cells.Item(15, 9).Value = 1.027681927861248
$ws.Cells.Item(15, 10).Value = 1.031315074403884
$ws.Cells.Item(15, 11).Value = 1.026882192518696
$ws.Cells.Item(15, 12).Value = 1.036926835544606
$ws.Cells.Item(15, 13).Value = 1.044098568512557
$ws.Cells.Item(15, 14).Value = 1.013796284965319

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.02603390762909
$ws.Cells.Item(16, 4).Value = 1.023634784689579
$ws.Cells.Item(16, 5).Value = 1.034584864946816
$ws.Cells.Item(16, 6).Value = 1.041910920141611
$ws.Cells.Item(16, 9).Value = 1.027746679784027
$ws.Cells.Item(16, 10).Value = 1.032419879352862
$ws.Cells.Item(16, 11).Value = 1.02711695452225
$ws.Cells.Item(16, 12).Value = 1.038027428624988
$ws.Cells.Item(16, 13).Value = 1.04532749814416
$ws.Cells.Item(16, 14).Value = 1.014192648055195

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.02686745018941
$ws.Cells.Item(17, 4).Value = 1.023853086440131
$ws.Cells.Item(17, 5).Value = 1.035347855050772
$ws.Cells.Item(17, 6).Value = 1.042754090270336
$ws.Cells.Item(17, 9).Value = 1.027785696793537
$ws.Cells.Item(17, 10).Value = 1.03311102144153
$ws.Cells.Item(17, 11).Value = 1.027261812832631
$ws.Cells.Item(17, 12).Value = 1.038715969682559
$ws.Cells.Item(17, 13).Value = 1.046096560209594
$ws.Cells.Item(17, 14).Value = 1.014440296640642

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.027352967761147
$ws.Cells.Item(18, 4).Value = 1.02397963304905
$ws.Cells.Item(18, 5).Value = 1.035792284708298
$ws.Cells.Item(18, 6).Value = 1.043245289960371
$ws.Cells.Item(18, 9).Value = 1.027807877329898
$ws.Cells.Item(18, 10).Value = 1.033513487843186
$ws.Cells.Item(18, 11).Value = 1.027345439387843
$ws.Cells.Item(18, 12).Value = 1.039116933547572
$ws.Cells.Item(18, 13).Value = 1.046544498368184
$ws.Cells.Item(18, 14).Value = 1.014584396272636

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027518403541067
$ws.Cells.Item(19, 4).Value = 1.024022648879112
$ws.Cells.Item(19, 5).Value = 1.035943721312831
$ws.Cells.Item(19, 6).Value = 1.043412674438442
$ws.Cells.Item(19, 9).Value = 1.027815342305293
$ws.Cells.Item(19, 10).Value = 1.033650606598615
$ws.Cells.Item(19, 11).Value = 1.027373806722036
$ws.Cells.Item(19, 12).Value = 1.039253542438254
$ws.Cells.Item(19, 13).Value = 1.04669712568265
$ws.Cells.Item(19, 14).Value = 1.014633471543508

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.026778088910315
$ws.Cells.Item(20, 4).Value = 1.023829745906674
$ws.Cells.Item(20, 5).Value = 1.035266056715309
$ws.Cells.Item(20, 6).Value = 1.042663689154085
$ws.Cells.Item(20, 9).Value = 1.027781570346447
$ws.Cells.Item(20, 10).Value = 1.033036937450378
$ws.Cells.Item(20, 11).Value = 1.027246360539083
$ws.Cells.Item(20, 12).Value = 1.038642163243583
$ws.Cells.Item(20, 13).Value = 1.046014113847933
$ws.Cells.Item(20, 14).Value = 1.014413762527972

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024363253327552
$ws.Cells.Item(21, 4).Value = 1.023193381061207
$ws.Cells.Item(21, 5).Value = 1.03305566662043
$ws.Cells.Item(21, 6).Value = 1.040221449297817
$ws.Cells.Item(21, 9).Value = 1.027664993056494
$ws.Cells.Item(21, 10).Value = 1.03103395356437
$ws.Cells.Item(21, 11).Value = 1.026821841561942
$ws.Cells.Item(21, 12).Value = 1.036646796958579
$ws.Cells.Item(21, 13).Value = 1.043785947385488
$ws.Cells.Item(21, 14).Value = 1.01369533448756

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.022838626403462
$ws.Cells.Item(22, 4).Value = 1.022786227006253
$ws.Cells.Item(22, 5).Value = 1.031660192195374
$ws.Cells.Item(22, 6).Value = 1.038680200795731
$ws.Cells.Item(22, 9).Value = 1.027586509990163
$ws.Cells.Item(22, 10).Value = 1.029768393883549
$ws.Cells.Item(22, 11).Value = 1.026547147649343
$ws.Cells.Item(22, 12).Value = 1.035386161784353
$ws.Cells.Item(22, 13).Value = 1.042378990933125
$ws.Cells.Item(22, 14).Value = 1.013240407459075

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.023647471882663
$ws.Cells.Item(23, 4).Value = 1.023002733729557
$ws.Cells.Item(23, 5).Value = 1.032400512166615
$ws.Cells.Item(23, 6).Value = 1.039497799546493
$ws.Cells.Item(23, 9).Value = 1.027628606217047
$ws.Cells.Item(23, 10).Value = 1.030439889504659
$ws.Cells.Item(23, 11).Value = 1.026693504457307
$ws.Cells.Item(23, 12).Value = 1.036055033862706
$ws.Cells.Item(23, 13).Value = 1.043125427002243
$ws.Cells.Item(23, 14).Value = 1.013481881958444

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.026818469506129
$ws.Cells.Item(24, 4).Value = 1.023840294915895
$ws.Cells.Item(24, 5).Value = 1.035303019744479
$ws.Cells.Item(24, 6).Value = 1.042704539402101
$ws.Cells.Item(24, 9).Value = 1.027783436694909
$ws.Cells.Item(24, 10).Value = 1.033070414871011
$ws.Cells.Item(24, 11).Value = 1.027253345444726
$ws.Cells.Item(24, 12).Value = 1.038675515205249
$ws.Cells.Item(24, 13).Value = 1.046051369794693
$ws.Cells.Item(24, 14).Value = 1.01442575323124

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030471933579147
$ws.Cells.Item(25, 4).Value = 1.024781422877931
$ws.Cells.Item(25, 5).Value = 1.03864741647204
$ws.Cells.Item(25, 6).Value = 1.046402085227728
$ws.Cells.Item(25, 9).Value = 1.027940447854658
$ws.Cells.Item(25, 10).Value = 1.03609700809811
$ws.Cells.Item(25, 11).Value = 1.027868987877461
$ws.Cells.Item(25, 12).Value = 1.041691023941352
$ws.Cells.Item(25, 13).Value = 1.049421663745579
$ws.Cells.Item(25, 14).Value = 1.015507385343153
